# New activity data (row 2) plus the pre-existing rows shifted down by one.
# Columns: A Name, B Distance(km), C Moving Time (seconds, numeric now), D Type,
#          E Start Date, F Avg Speed (min/km), G Avg Speed (m/s), H Avg HR (b/min)
$data = @(
    @("Intervallit",      11.05, 3267, "Run", "2025-07-10T11:46:41Z", "04:56", 3.383, 158.5),
    @("Lunch Run",        12,    4688, "Run", "2025-07-09T12:24:49Z", "06:31", 2.561, 139.5),
    @("Säbä",             6.41,  5900, "Run", "2025-07-08T19:33:57Z", "15:20", 1.086, 131.2),
    @("Night Run",        9.09,  3666, "Run", "2025-07-07T21:57:25Z", "06:43", 2.479, 138.6),
    @("Afternoon Run",    21.27, 7807, "Run", "2025-07-05T16:17:54Z", "06:07", 2.724, 145.7),
    @("Vetoja",           7.06,  2404, "Run", "2025-07-03T20:37:43Z", "05:41", 2.938, 153.8),
    @("Morning Run",      10.65, 4308, "Run", "2025-07-02T10:27:49Z", "06:45", 2.471, 135),
    @("Säbä",             5.03,  4714, "Run", "2025-07-01T20:42:49Z", "15:37", 1.067, 142.1),
    @("Alkulämpöfutis",   1.4,   760,  "Run", "2025-07-01T20:24:46Z", "09:03", 1.839, 135.1)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
